{"js": "// Office.js (Word JavaScript API) script.\n// Applies proofing-mark (spellcheck/grammar) run-splits and trailing-period\n// fixes to the \"Hints\" list, matching the target OOXML diff.\n//\n// Strategy: each touched paragraph is rewritten in place (insertOoxml with\n// location \"Replace\") using a minimal, namespaced <w:p> fragment that\n// reproduces the exact <w:r>/<w:proofErr> run layout from the diff. This\n// lets us emit <w:proofErr .../> markers that have no first-class Word\n// JS API surface, while leaving every other paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Namespace declarations Word expects on the root of an OOXML fragment.\nconst WNS =\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapOoxml(pXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document ' + WNS + \"><w:body>\" +\n    pXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// text -> exact replacement <w:p>...</w:p> inner-run XML (no w:p wrapper).\nconst replacements = [\n  [\n    \"Triforce: Remember Layering (and Princess Zelda)\",\n    '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>Triforce</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"<w:r><w:t>: Remember Layering (and Princess Zelda)</w:t></w:r>\",\n  ],\n  [\n    \"Can You Hear Me Now: Good... Now use Rectangles.\",\n    \"<w:r><w:t>Can You Hear Me Now: Good</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>...</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> Now use Rectangles.</w:t></w:r>',\n  ],\n  [\n    \"Native American Example 8: A great example of Triangle Iteration\",\n    \"<w:r><w:t>Native American Ex</w:t></w:r>\" +\n      \"<w:r><w:t>ample 8</w:t></w:r>\" +\n      \"<w:r><w:t>: A great example of Triangle Iteration</w:t></w:r>\" +\n      \"<w:r><w:t>.</w:t></w:r>\",\n  ],\n  [\n    \"Tie Fighter: Darth Vader's Tie Fighter was special because it's wings were made with Linear Iteration.\",\n    '<w:r><w:t xml:space=\"preserve\">Tie Fighter: Darth Vader\\u2019s Tie Fighter was special because </w:t></w:r>'\n      .replace(\"\\u2019\", \"'\") +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>it's</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> wings were made with Linear Iteration.</w:t></w:r>',\n  ],\n  [\n    \"Moon: Use Circles to fight evil by moonlight\",\n    \"<w:r><w:t>Moon: Use Circles to fight evil by moonlight</w:t></w:r>\" +\n      \"<w:r><w:t>.</w:t></w:r>\",\n  ],\n  [\n    \"ROSS: It's ROSS the bunny!  Bunny Adventure coming soon to iPhone!\",\n    '<w:r><w:t xml:space=\"preserve\">ROSS: It\\'s ROSS the bunny!  Bunny Adventure coming soon to </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>iPhone</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"<w:r><w:t>!</w:t></w:r>\",\n  ],\n  [\n    \"Native American Example 3: Do not make your background too big and use some creative iterations.\",\n    '<w:r><w:t xml:space=\"preserve\">Native American Example 3: Do not make your background too big and use some creative </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>iterations</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"<w:r><w:t>.</w:t></w:r>\",\n  ],\n  [\n    \"Megaman: He is a super fighting robot\",\n    '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>Megaman</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"<w:r><w:t>: He is a super fighting robot</w:t></w:r>\" +\n      \"<w:r><w:t>.</w:t></w:r>\",\n  ],\n];\n\nfor (const para of paragraphs.items) {\n  for (const [needle, innerXml] of replacements) {\n    if (para.text === needle) {\n      para.insertOoxml(wrapOoxml(\"<w:p>\" + innerXml + \"</w:p>\"), \"Replace\");\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies proofing-mark (spellcheck/grammar) run-splits and trailing-period\n# fixes to the \"Hints\" list, matching the target OOXML diff.\n#\n# Strategy: each touched paragraph's Range is rewritten in place with\n# Range.InsertXML, feeding a minimal <pkg:package>/<w:document> fragment\n# that reproduces the exact <w:r>/<w:proofErr> run layout from the diff.\n# InsertXML replaces the addressed range's content, so the paragraph\n# count and surrounding paragraphs are left untouched.\n\n$d = $word.ActiveDocument\n\nfunction New-OoxmlFragment([string]$innerParagraphXml) {\n    $ns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n    return '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        \"<w:document $ns><w:body><w:p>$innerParagraphXml</w:p></w:body></w:document>\" +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$replacements = @{\n    \"Triforce: Remember Layering (and Princess Zelda)\" =\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Triforce</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>: Remember Layering (and Princess Zelda)</w:t></w:r>'\n\n    \"Can You Hear Me Now: Good... Now use Rectangles.\" =\n        '<w:r><w:t>Can You Hear Me Now: Good</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>...</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> Now use Rectangles.</w:t></w:r>'\n\n    \"Native American Example 8: A great example of Triangle Iteration\" =\n        '<w:r><w:t>Native American Ex</w:t></w:r>' +\n        '<w:r><w:t>ample 8</w:t></w:r>' +\n        '<w:r><w:t>: A great example of Triangle Iteration</w:t></w:r>' +\n        '<w:r><w:t>.</w:t></w:r>'\n\n    \"Tie Fighter: Darth Vader's Tie Fighter was special because it's wings were made with Linear Iteration.\" =\n        '<w:r><w:t xml:space=\"preserve\">Tie Fighter: Darth Vader''s Tie Fighter was special because </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>it''s</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> wings were made with Linear Iteration.</w:t></w:r>'\n\n    \"Moon: Use Circles to fight evil by moonlight\" =\n        '<w:r><w:t>Moon: Use Circles to fight evil by moonlight</w:t></w:r>' +\n        '<w:r><w:t>.</w:t></w:r>'\n\n    \"ROSS: It's ROSS the bunny!  Bunny Adventure coming soon to iPhone!\" =\n        '<w:r><w:t xml:space=\"preserve\">ROSS: It''s ROSS the bunny!  Bunny Adventure coming soon to </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>iPhone</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>!</w:t></w:r>'\n\n    \"Native American Example 3: Do not make your background too big and use some creative iterations.\" =\n        '<w:r><w:t xml:space=\"preserve\">Native American Example 3: Do not make your background too big and use some creative </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>iterations</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:t>.</w:t></w:r>'\n\n    \"Megaman: He is a super fighting robot\" =\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Megaman</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>: He is a super fighting robot</w:t></w:r>' +\n        '<w:r><w:t>.</w:t></w:r>'\n}\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($text)) {\n        $xml = New-OoxmlFragment $replacements[$text]\n        [void]$p.Range.InsertXML($xml)\n    }\n}\n"}
